$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header changes
$ws.Range("C1").Value = "rules"
$ws.Range("E1").Value = "adaptive_filter"

# Row 2
$ws.Range("E2").Value = "RLS"
$ws.Range("F2").Value = [double]"1.136128162321535e+140"
$ws.Range("G2").Value = [double]"1.896339335057897e+139"
$ws.Range("H2").Value = [double]"1.13600797758886e+140"

# Row 3
$ws.Range("E3").Value = "RLS"
$ws.Range("F3").Value = 505.0920875197271
$ws.Range("G3").Value = 84.30615710053063
$ws.Range("H3").Value = 481.8763475666818

# Row 4
$ws.Range("E4").Value = "RLS"
$ws.Range("F4").Value = 6.556937094707963
$ws.Range("G4").Value = 1.094434425847458
$ws.Range("H4").Value = 5.365066719927139

# Row 5
$ws.Range("E5").Value = "RLS"
$ws.Range("F5").Value = 6.213408710176207
$ws.Range("G5").Value = 1.037095262018853
$ws.Range("H5").Value = 4.926389891355839

# Row 6
$ws.Range("E6").Value = "RLS"
$ws.Range("F6").Value = 5.987407937800259
$ws.Range("G6").Value = 0.9993729197142451
$ws.Range("H6").Value = 4.584053117668564

# Row 7
$ws.Range("E7").Value = "RLS"
$ws.Range("F7").Value = 6.025546986257404
$ws.Range("G7").Value = 1.005738801679813
$ws.Range("H7").Value = 4.798328146692079
